$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Datos actualizados" timestamp string (A1) ---
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 25 de Agosto de 2020 a las 16:25"

# --- Country name swaps caused by shared-string reshuffle in source diff ---
# Namibia / Zimbabue swap (rows 105/106)
$ws.Cells.Item(105, 1).Value = "Namibia"
$ws.Cells.Item(106, 1).Value = "Zimbabue"

# Gambia / Eslovenia swap (rows 130/131)
$ws.Cells.Item(130, 1).Value = "Gambia"
$ws.Cells.Item(131, 1).Value = "Eslovenia"

# Montserrat / Islas Malvinas swap (rows 214/215)
$ws.Cells.Item(214, 1).Value = "Montserrat"
$ws.Cells.Item(215, 1).Value = "Islas Malvinas"

# --- Updated COVID-19 stat counters (refreshed data pull) ---
# Row 4
$ws.Cells.Item(4, 2).Value = 5918213
$ws.Cells.Item(4, 3).Value = 2583
$ws.Cells.Item(4, 4).Value = 3219327
$ws.Cells.Item(4, 5).Value = 2517660
$ws.Cells.Item(4, 7).Value = 112
$ws.Cells.Item(4, 8).Value = 181226

# Row 6
$ws.Cells.Item(6, 2).Value = 3191977
$ws.Cells.Item(6, 3).Value = 27096
$ws.Cells.Item(6, 4).Value = 2425641
$ws.Cells.Item(6, 5).Value = 707543
$ws.Cells.Item(6, 7).Value = 247
$ws.Cells.Item(6, 8).Value = 58793

# Row 24
$ws.Cells.Item(24, 2).Value = 211947
$ws.Cells.Item(24, 3).Value = 3962
$ws.Cells.Item(24, 4).Value = 153761
$ws.Cells.Item(24, 5).Value = 51590
$ws.Cells.Item(24, 7).Value = 77
$ws.Cells.Item(24, 8).Value = 6596

# Row 33
$ws.Cells.Item(33, 5).Value = 9974
$ws.Cells.Item(33, 8).Value = 1523

# Row 49
$ws.Cells.Item(49, 4).Value = 54816
$ws.Cells.Item(49, 5).Value = 1592

# Row 50
$ws.Cells.Item(50, 2).Value = 55912
$ws.Cells.Item(50, 3).Value = 192
$ws.Cells.Item(50, 4).Value = 41021
$ws.Cells.Item(50, 5).Value = 13086
$ws.Cells.Item(50, 7).Value = 4
$ws.Cells.Item(50, 8).Value = 1805

# Row 61
$ws.Cells.Item(61, 5).Value = 3860
$ws.Cells.Item(61, 7).Value = 1
$ws.Cells.Item(61, 8).Value = 2002

# Row 66
$ws.Cells.Item(66, 2).Value = 34358
$ws.Cells.Item(66, 3).Value = 530
$ws.Cells.Item(66, 4).Value = 23869
$ws.Cells.Item(66, 5).Value = 9529
$ws.Cells.Item(66, 7).Value = 15
$ws.Cells.Item(66, 8).Value = 960

# Row 99
$ws.Cells.Item(99, 2).Value = 8379
$ws.Cells.Item(99, 3).Value = 33
$ws.Cells.Item(99, 4).Value = 7179
$ws.Cells.Item(99, 5).Value = 1133

# Row 105
$ws.Cells.Item(105, 2).Value = 6160
$ws.Cells.Item(105, 3).Value = 130
$ws.Cells.Item(105, 4).Value = 2732
$ws.Cells.Item(105, 5).Value = 3371
$ws.Cells.Item(105, 7).Value = 1
$ws.Cells.Item(105, 8).Value = 57

# Row 106
$ws.Cells.Item(106, 2).Value = 6070
$ws.Cells.Item(106, 4).Value = 4950
$ws.Cells.Item(106, 5).Value = 965
$ws.Cells.Item(106, 8).Value = 155

# Row 130
$ws.Cells.Item(130, 3).Value = 101
$ws.Cells.Item(130, 4).Value = 601
$ws.Cells.Item(130, 5).Value = 1995
$ws.Cells.Item(130, 7).Value = 3
$ws.Cells.Item(130, 8).Value = 90

# Row 131
$ws.Cells.Item(131, 2).Value = 2686
$ws.Cells.Item(131, 3).Value = 21
$ws.Cells.Item(131, 4).Value = 2139
$ws.Cells.Item(131, 5).Value = 414
$ws.Cells.Item(131, 8).Value = 133

# Row 214
$ws.Cells.Item(214, 4).Value = 12
$ws.Cells.Item(214, 8).Value = 1

# Row 215
$ws.Cells.Item(215, 4).Value = 13
$ws.Cells.Item(215, 8).Value = 0
